$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The RxCUI code for Flurbiprofen (row 2) was recorded incorrectly as the
# PharmGKB-style string "C0016377". Correct it to the actual numeric RxCUI.
$ws.Range("D2").Value = 4502

# Update the active selection to match the saved workbook state.
[void]$ws.Range("B11").Select()
